# glossary and comms data update @PC
# Add newly-recoded trust fund rows to the Trustee_Names glossary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trustee_Names")

# Carry the existing data-row style (used by rows 2:21) down into the
# new rows 22:29 so the appended records look like the rest of the table
# instead of the old (now-unused) blank-row formatting.
$ws.Range("A21:B21").Copy()
$ws.Range("A22:B29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 22
$ws.Range("A22").Value = "TF070611"
$ws.Range("B22").Value = "MDTF"

# Row 23 (string discovery order: description text before the TF code)
$ws.Range("B23").Value = "Standby Recovery Financing Facility"
$ws.Range("A23").Value = "TF070948"

# Row 24
$ws.Range("A24").Value = "TF070809"
$ws.Range("B24").Value = "Japan TF for Mainstreaming DRM"

# Row 25
$ws.Range("A25").Value = "TF070952"
$ws.Range("B25").Value = "South-South Cooperation MDTF"

# Row 26
$ws.Range("A26").Value = "TF070806"
$ws.Range("B26").Value = "Spanish Trust Fund for Mainstreaming DRR"

# Row 27
$ws.Range("A27").Value = "TF070868"
$ws.Range("B27").Value = "Callable Funds - Standby Recovery Financing Facility "

# Row 28
$ws.Range("A28").Value = "TF071345"
$ws.Range("B28").Value = "Track III Standby Recovery Financing Facility"

# Row 29
$ws.Range("A29").Value = "TF070807"
$ws.Range("B29").Value = "Australian Trust Fund for Mainstreaming DRR"

# Match the author's final selection/cursor position.
$null = $ws.Range("B30").Select()
